# Revert the payslip data sheet back to the "local development" dataset.
# Replaces the content of rows 2-8 (employees) with a different set of
# records, drops the former row 9 entirely (dimension shrinks to A1:M8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Remove the last data row (old row 9) completely - the new dataset only
# has 7 employee rows (rows 2-8) instead of 8 (rows 2-9).
# ---------------------------------------------------------------------
$ws.Rows(9).Delete()

# ---------------------------------------------------------------------
# Helper-free, per-cell assignment of the new values. Columns are:
#   A ID, B FirstName, C LastName, D Rate, E Days Worked, F Total Amount,
#   G Mess, H Advance, I Home Advance, J Sunday Expenditure, K Net Pay,
#   L Creation Date
# ---------------------------------------------------------------------

# Row 2 - SHIV SHANKAR MANJHI
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "SHIV SHANKAR "
$ws.Cells.Item(2,3).Value = "MANJHI"
$ws.Cells.Item(2,4).Value = 600
$ws.Cells.Item(2,5).Value = 23
$ws.Cells.Item(2,6).Value = 13800
$ws.Cells.Item(2,7).Value = 2400.0
$ws.Cells.Item(2,8).Value = 5000.0
$ws.Cells.Item(2,9).Value = 3000.0
$ws.Cells.Item(2,10).Value = 1200.0
$ws.Cells.Item(2,11).Value = 2600
$ws.Cells.Item(2,12).Value = "2024-03-04 08:58:58"

# Row 3 - ROBIN MANDAL (Rate and Total Amount are blank in the new data)
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "ROBIN "
$ws.Cells.Item(3,3).Value = "MANDAL"
$ws.Cells.Item(3,4).ClearContents()
$ws.Cells.Item(3,5).Value = 30
$ws.Cells.Item(3,6).ClearContents()
$ws.Cells.Item(3,7).Value = 1200.0
$ws.Cells.Item(3,8).Value = 5000.0
$ws.Cells.Item(3,9).Value = 2000.0
$ws.Cells.Item(3,10).Value = 2400.0
$ws.Cells.Item(3,11).Value = -10600
$ws.Cells.Item(3,12).Value = "2024-03-04 09:04:11"

# Row 4 - KALI TUDU
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "KALI "
$ws.Cells.Item(4,3).Value = "TUDU"
$ws.Cells.Item(4,4).Value = 700
$ws.Cells.Item(4,5).Value = -10
$ws.Cells.Item(4,6).Value = 7000
$ws.Cells.Item(4,7).Value = 21.0
$ws.Cells.Item(4,8).Value = 211.0
$ws.Cells.Item(4,9).Value = 21.0
$ws.Cells.Item(4,10).Value = 12.0
$ws.Cells.Item(4,11).Value = 6851.6666666667
$ws.Cells.Item(4,12).Value = "2024-03-14 14:54:39"

# Row 5 - NEPAL MAHTO
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "NEPAL "
$ws.Cells.Item(5,3).Value = "MAHTO"
$ws.Cells.Item(5,4).Value = 800
$ws.Cells.Item(5,5).Value = 30
$ws.Cells.Item(5,6).Value = 24000
$ws.Cells.Item(5,7).Value = 2400.0
$ws.Cells.Item(5,8).Value = 1000.0
$ws.Cells.Item(5,9).Value = 4000.0
$ws.Cells.Item(5,10).Value = 1200.0
$ws.Cells.Item(5,11).Value = 15933.333333333
$ws.Cells.Item(5,12).Value = "2024-03-14 15:51:53"

# Row 6 - BANSHI MAHTO
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "BANSHI "
$ws.Cells.Item(6,3).Value = "MAHTO"
$ws.Cells.Item(6,4).Value = 700
$ws.Cells.Item(6,5).Value = 30
$ws.Cells.Item(6,6).Value = 21000
$ws.Cells.Item(6,7).Value = 3200.0
$ws.Cells.Item(6,8).Value = 1000.0
$ws.Cells.Item(6,9).Value = 1000.0
$ws.Cells.Item(6,10).Value = 1200.0
$ws.Cells.Item(6,11).Value = 12566.666666667
$ws.Cells.Item(6,12).Value = "2024-03-14 16:00:15"

# Row 7 - RAJESH ROSHAN (Rate is the literal text "600 ", not a number)
$ws.Cells.Item(7,1).Value = 8
$ws.Cells.Item(7,2).Value = "RAJESH "
$ws.Cells.Item(7,3).Value = "ROSHAN"
$ws.Cells.Item(7,4).Formula = "'600 "
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = 25
$ws.Cells.Item(7,6).Value = 15000.0
$ws.Cells.Item(7,7).Value = 1200.0
$ws.Cells.Item(7,8).Value = 1000.0
$ws.Cells.Item(7,9).Value = 7000.0
$ws.Cells.Item(7,10).Value = 1000.0
$ws.Cells.Item(7,11).Value = 5090.0
$ws.Cells.Item(7,12).Value = "2024-03-20 22:28:35"

# Row 8 - nassour h
$ws.Cells.Item(8,1).Value = 13
$ws.Cells.Item(8,2).Value = "nassour "
$ws.Cells.Item(8,3).Value = "h"
$ws.Cells.Item(8,4).Value = 123
$ws.Cells.Item(8,5).Value = 20
$ws.Cells.Item(8,6).Value = 2460
$ws.Cells.Item(8,7).Value = 2.0
$ws.Cells.Item(8,8).Value = 2.0
$ws.Cells.Item(8,9).Value = 2.0
$ws.Cells.Item(8,10).Value = 2.0
$ws.Cells.Item(8,11).Value = 2472.5
$ws.Cells.Item(8,12).Value = "2024-03-24 18:09:58"
